$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.434.71'
$ws.Range("E2").Value = '  +0.52%  '

$ws.Range("D3").Value = '1.908.02'
$ws.Range("E3").Value = '  -0.11%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.008'
$ws.Range("E4").Value = '  +0.69%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.98'

$ws.Range("E6").Value = '  +0.55%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4810'
$ws.Range("E7").Value = '  +1.85%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4057'
$ws.Range("E8").Value = '  -0.15%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08143'
$ws.Range("E9").Value = '  +1.43%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.012'
$ws.Range("E10").Value = '  +1.05%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '23.37'
$ws.Range("E11").Value = '  +3.87%  '

$ws.Range("D12").Value = '1.937.57'
$ws.Range("E12").Value = '  -0.65%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.998'
$ws.Range("E13").Value = '  +1.97%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.150'
$ws.Range("E14").Value = '  +0.41%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '90.18'
$ws.Range("E15").Value = '  +0.69%  '

$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.008'
$ws.Range("E16").Value = '  +0.70%  '

$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.06780'
$ws.Range("E17").Value = '  +2.21%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001036'
$ws.Range("E18").Value = '  +0.66%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.66'
$ws.Range("E19").Value = '  +0.24%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.005'
$ws.Range("E20").Value = '  +0.40%  '

$ws.Range("D21").Value = '29.457.41'
$ws.Range("E21").Value = '  +0.50%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.611'
$ws.Range("E22").Value = '  +1.60%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.72'
$ws.Range("E23").Value = '  +2.29%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.181'
$ws.Range("E24").Value = '  -0.78%  '

$ws.Range("D25").Value = '2.143.11'
$ws.Range("E25").Value = '  -1.68%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '155.59'
$ws.Range("E26").Value = '  +0.65%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.409'
$ws.Range("E27").Value = '  +6.71%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.02'
$ws.Range("E28").Value = '  +1.11%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.106'
$ws.Range("E29").Value = '  +0.04%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '119.98'
$ws.Range("E30").Value = '  +2.19%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.023'
$ws.Range("E31").Value = '  -4.03%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09527'
$ws.Range("E32").Value = '  +0.37%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.528'
$ws.Range("E33").Value = '  +2.79%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.558'
$ws.Range("E34").Value = '  +0.53%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.386'
$ws.Range("E35").Value = '  -2.63%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02266'
$ws.Range("E36").Value = '  +0.86%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06089'
$ws.Range("E37").Value = '  +0.23%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.175'
$ws.Range("E38").Value = '  +0.36%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '10.81'
$ws.Range("E39").Value = '  +7.14%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5951'
$ws.Range("E40").Value = '  +1.51%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.987'
$ws.Range("E41").Value = '  -2.91%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1851'
$ws.Range("E42").Value = '  +0.89%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.277'
$ws.Range("E43").Value = '  +0.30%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.371'
$ws.Range("E44").Value = '  -6.04%  '

$ws.Range("E45").Value = '  +3.81%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.07590'
$ws.Range("E46").Value = '  -3.45%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5564'
$ws.Range("E47").Value = '  +0.76%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.941'
$ws.Range("E48").Value = '  +1.19%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '115.92'
$ws.Range("E49").Value = '  +2.52%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '72.46'
$ws.Range("E50").Value = '  +1.62%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.399'
$ws.Range("E51").Value = '  +2.23%  '
